$d = $word.ActiveDocument

# Each Sprint Retrospective question paragraph currently ends in a bare
# ":" run (italic, Arial) right after a hyperlinked reviewer name. The
# edit appends each reviewer's comment text directly after that colon,
# inside the very same run (so its formatting - rFonts Arial, italic,
# rtl=0 - is preserved exactly as authored).
#
# Trick: Word.Range.Text assignment takes on the formatting of the
# *character to its left* when the range is collapsed (zero-length) and
# positioned strictly inside an existing run. Using Paragraph.Range.End-1
# (i.e. just before the paragraph mark, right after the ":") keeps the
# insertion welded onto the colon's own run instead of bleeding in the
# neighbouring hyperlink's formatting.

$edits = @(
    @{ Match = "Samson Cournane:"; Append = "  I think that our team did a good job at collaborating and help us finish  our goal on time. Our team very good at working together with our PM." },
    @{ Match = "Sam Minor:"; Append = " I think the team has done a good job in working together to meet our goals, looking towards the deliverables. " },
    @{ Match = "Tereza Holubcova:"; Append = " The team has been working very well together, thanks to our PM. He has been such a great leader and I bet he did a lot of work on its own. I feel like the team is pretty dedicated, open and easy to work with, which makes the workload more manageable." },
    @{ Match = "Emily Brule:"; Append = " I think that our ability to communicate together as a team is great and really helpful for our project. The workload and directions can be tricky at times but our ability to work well together is like our saving grace." },
    @{ Match = "Samson Cournane:"; Append = " I think that sometimes when we work together online it can be a bit challenging, but I think that if we work together in-person collaboration would be easier." },
    @{ Match = "Sam Minor:"; Append = " Sometimes the way we communicate in a team setting can be a bit ambigous" },
    @{ Match = "Tereza Holubcova:"; Append = " I agree with Andrew about knowing what is “good enough” or if it is meeting the assignment requirements." },
    @{ Match = "Emily Brule:"; Append = " I think that collaborating without really any opportunities to have class together or meet in person makes working on such a difficult project even more difficult." }
)

$editIndex = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    if ($editIndex -ge $edits.Count) {
        break
    }
    $p = $d.Paragraphs.Item($i)
    $ptext = $p.Range.Text.TrimEnd([char]13)
    $want = $edits[$editIndex].Match
    if ($ptext -eq $want) {
        $endPos = $p.Range.End - 1
        $r = $d.Range($endPos, $endPos)
        $r.Text = $edits[$editIndex].Append
        $editIndex = $editIndex + 1
    }
}

Write-Host ("Applied " + $editIndex + " of " + $edits.Count + " edits")
